$wb = $excel.ActiveWorkbook
$wb.Worksheets.Item("profesores").Activate()
$ws = $wb.ActiveSheet

# New "Resultados finales" numbers (ID / ejecución) added across the
# professors sheet. Each pair is (row, B-value, C-value); a blank string
# means "leave that column untouched".

$ws.Range("B2").Value = 51837750
$ws.Range("C2").Value = 51837750

$ws.Range("B3").Value = 52857739
$ws.Range("C3").Value = 406527

$ws.Range("C7").Value = 558959

$ws.Range("C9").Value = 79885769

$ws.Range("B13").Value = 1024479509
$ws.Range("C13").Value = 406593

$ws.Range("B14").Value = 41798596
$ws.Range("C14").Value = 406544

$ws.Range("C15").Value = 79403099

$ws.Range("C17").Value = 801363

$ws.Range("C18").Value = 27078

$ws.Range("B20").Value = 39762772

$ws.Range("B21").Value = 1022336186

$ws.Range("B22").Value = 1018423989

# Final selection left on B22, matching where data entry ended.
$ws.Range("B22").Select()
